$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Flight Page"
$ws.Range("B3").Value = 10.444
